$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the worksheet to reflect its content
$ws.Name = "Keycap Dimensions"

# --- Physical verification pass against real switches/keycaps ---

# k_width: standard keycap width
$ws.Range("B1").Value = 18.2

# k_poledia: keycap pole diameter
$ws.Range("B11").Value = 5.3

# k_cmxw: keycap pole Cherry MX width (formula referencing clearance cell)
$ws.Range("B12").Formula = "=1.3 + `$I`$16"

# k_cmxd: keycap Cherry MX cut-in depth - reduced for extra strength
$ws.Range("B14").Value = 4

# s_bwidth: keyswitch base width (inconsequential)
$ws.Range("B16").Value = 14.4

# s_polebw: keyswitch pole base width
$ws.Range("B18").Value = 5.4

# s_polebl: keyswitch pole base length (was mislabeled "Height")
$ws.Range("B19").Value = 7.45
$ws.Range("D19").Value = "Keyswitch Pole Base Length"

# Leave the selection where the author last left it
$ws.Range("B12").Select()
